# Actualización automática 2025-07-21 14:45:09
# Rotate advisor data for rows 7-9 on "VENTAS POR GRUPO" and "VENTA MENSUAL".

$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO (only names rotate; numeric columns are all 0) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("B7").Value = "LOZANO MOLINA TITO JERSON"
$ws1.Range("B8").Value = "MACHARE BARCO LISSETTE STEFANIA"
$ws1.Range("B9").Value = "MEZA FERNANDEZ JONATHAN ALEXIS"

# --- Sheet: VENTA MENSUAL (names + monthly values rotate together) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("B7").Value = "LOZANO MOLINA TITO JERSON"
$ws2.Range("C7").Value = 144.53
$ws2.Range("D7").Value = 11.52
$ws2.Range("E7").Value = 10.44

$ws2.Range("B8").Value = "MACHARE BARCO LISSETTE STEFANIA"
$ws2.Range("C8").Value = 0
$ws2.Range("D8").Value = 178.33
$ws2.Range("E8").Value = 0

$ws2.Range("B9").Value = "MEZA FERNANDEZ JONATHAN ALEXIS"
$ws2.Range("C9").Value = 0
$ws2.Range("D9").Value = 0
$ws2.Range("E9").Value = 0
